# Update market-price figures (currentAveragePrice / LevePrice / LeveProfit columns)
# for the Seraph_Profits leve-crafting sheets, per the latest scheduled market-board pull.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")

# Row 9
$ws.Range("H9").Value = 230
$ws.Range("I9").Value = 230
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 230
$ws.Range("L9").Value = 0
$ws.Range("M9").Value = -61
$ws.Range("N9").ClearContents()

# Row 64
$ws.Range("H64").Value = 3360
$ws.Range("I64").Value = 4000

# Row 67
$ws.Range("H67").Value = 3360
$ws.Range("I67").Value = 4000

# Row 137
$ws.Range("H137").Value = 1238.9286
$ws.Range("I137").Value = 1156.5
$ws.Range("J137").Value = 1445
$ws.Range("K137").Value = 3469.5
$ws.Range("L137").Value = 4335
$ws.Range("M137").Value = -919.5
$ws.Range("N137").Value = -9435

$ws = $wb.Worksheets.Item("ARM")

# Row 2
$ws.Range("H2").Value = 4896.5
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 4896.5
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 4896.5
$ws.Range("M2").ClearContents()
$ws.Range("N2").Value = -5122.5

# Row 74
$ws.Range("H74").Value = 704.3103599999999
$ws.Range("I74").Value = 719.62964
$ws.Range("J74").Value = 497.5
$ws.Range("K74").Value = 719.62964
$ws.Range("L74").Value = 497.5
$ws.Range("M74").Value = 154.37036
$ws.Range("N74").Value = -2245.5

# Row 77
$ws.Range("H77").Value = 704.3103599999999
$ws.Range("I77").Value = 719.62964
$ws.Range("J77").Value = 497.5
$ws.Range("K77").Value = 3598.1482
$ws.Range("L77").Value = 2487.5
$ws.Range("M77").Value = 769.8517999999999
$ws.Range("N77").Value = -11223.5

# Row 102
$ws.Range("H102").Value = 2113.4
$ws.Range("I102").Value = 2022.6666
$ws.Range("J102").Value = 2249.5
$ws.Range("K102").Value = 2022.6666
$ws.Range("L102").Value = 2249.5
$ws.Range("M102").Value = -400.6666
$ws.Range("N102").Value = -5493.5

# Row 116
$ws.Range("H116").Value = 4896.5
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 4896.5
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 4896.5
$ws.Range("M116").ClearContents()
$ws.Range("N116").Value = -9484.5

# Row 132
$ws.Range("H132").Value = 20100.545
$ws.Range("I132").Value = 2345.111
$ws.Range("J132").Value = 100000
$ws.Range("K132").Value = 7035.333
$ws.Range("L132").Value = 300000
$ws.Range("M132").Value = -4505.333
$ws.Range("N132").Value = -305060

$ws = $wb.Worksheets.Item("BSM")

# Row 3
$ws.Range("H3").Value = 4896.5
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 4896.5
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 4896.5
$ws.Range("M3").ClearContents()
$ws.Range("N3").Value = -5124.5

# Row 134
$ws.Range("H134").Value = 2031.7693
$ws.Range("I134").Value = 2110.6667
$ws.Range("J134").Value = 1964.1428
$ws.Range("K134").Value = 6332.000100000001
$ws.Range("L134").Value = 5892.428400000001
$ws.Range("M134").Value = -3797.000100000001
$ws.Range("N134").Value = -10962.4284

$ws = $wb.Worksheets.Item("CRP")

# Row 31
$ws.Range("H31").Value = 6663.3335
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 6663.3335
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 6663.3335
$ws.Range("M31").ClearContents()
$ws.Range("N31").Value = -7253.3335

# Row 34
$ws.Range("H34").Value = 6663.3335
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 6663.3335
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 6663.3335
$ws.Range("M34").ClearContents()
$ws.Range("N34").Value = -7067.3335

# Row 58
$ws.Range("H58").Value = 2327.48
$ws.Range("I58").Value = 1317.7222
$ws.Range("J58").Value = 4924
$ws.Range("K58").Value = 1317.7222
$ws.Range("L58").Value = 4924
$ws.Range("M58").Value = -1114.7222
$ws.Range("N58").Value = -5330

# Row 105
$ws.Range("H105").Value = 4169.0967
$ws.Range("I105").Value = 716.26666
$ws.Range("J105").Value = 7406.125
$ws.Range("K105").Value = 716.26666
$ws.Range("L105").Value = 7406.125
$ws.Range("M105").Value = 1030.73334
$ws.Range("N105").Value = -10900.125

# Row 125
$ws.Range("H125").Value = 47499
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 47499
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 47499
$ws.Range("N125").Value = -52419

# Row 134
$ws.Range("H134").Value = 2726.5264
$ws.Range("I134").Value = 2616.4443
$ws.Range("J134").Value = 2996.7273
$ws.Range("K134").Value = 7849.3329
$ws.Range("L134").Value = 8990.1819
$ws.Range("M134").Value = -5314.3329
$ws.Range("N134").Value = -14060.1819

# Row 136
$ws.Range("H136").Value = 2327.48
$ws.Range("I136").Value = 1317.7222
$ws.Range("J136").Value = 4924
$ws.Range("K136").Value = 3953.1666
$ws.Range("L136").Value = 14772
$ws.Range("M136").Value = -1403.1666
$ws.Range("N136").Value = -19872

$ws = $wb.Worksheets.Item("CUL")

# Row 59
$ws.Range("H59").Value = 5000
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 5000
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 15000
$ws.Range("M59").ClearContents()
$ws.Range("N59").Value = -16080

# Row 107
$ws.Range("H107").Value = 143645.14
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 143645.14
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 430935.42
$ws.Range("N107").Value = -434775.42

$ws = $wb.Worksheets.Item("GSM")

# Row 132
$ws.Range("H132").Value = 3144.5715
$ws.Range("I132").Value = 2807.4
$ws.Range("J132").Value = 3987.5
$ws.Range("K132").Value = 8422.200000000001
$ws.Range("L132").Value = 11962.5
$ws.Range("M132").Value = -5892.200000000001
$ws.Range("N132").Value = -17022.5

$ws = $wb.Worksheets.Item("LTW")

# Row 55
$ws.Range("H55").Value = 511.1111
$ws.Range("I55").Value = 198.66667
$ws.Range("J55").Value = 1136
$ws.Range("K55").Value = 198.66667
$ws.Range("L55").Value = 1136
$ws.Range("M55").Value = -25.66667000000001
$ws.Range("N55").Value = -1482

# Row 61
$ws.Range("H61").Value = 2555.875
$ws.Range("I61").Value = 2421.0715
$ws.Range("J61").Value = 3499.5
$ws.Range("K61").Value = 2421.0715
$ws.Range("L61").Value = 3499.5
$ws.Range("M61").Value = -2219.0715
$ws.Range("N61").Value = -3903.5

# Row 62
$ws.Range("H62").Value = 17999
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 17999
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 17999
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -19247

# Row 65
$ws.Range("H65").Value = 17999
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 17999
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 53997
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -60237

# Row 93
$ws.Range("H93").Value = 1076.7646
$ws.Range("I93").Value = 893.6667
$ws.Range("J93").Value = 1516.2
$ws.Range("K93").Value = 893.6667
$ws.Range("L93").Value = 1516.2
$ws.Range("M93").Value = 354.3333
$ws.Range("N93").Value = -4012.2

# Row 106
$ws.Range("H106").Value = 17597
$ws.Range("I106").Value = 0
$ws.Range("J106").Value = 17597
$ws.Range("K106").Value = 0
$ws.Range("L106").Value = 17597
$ws.Range("N106").Value = -20121

# Row 113
$ws.Range("H113").Value = 2555.875
$ws.Range("I113").Value = 2421.0715
$ws.Range("J113").Value = 3499.5
$ws.Range("K113").Value = 2421.0715
$ws.Range("L113").Value = 3499.5
$ws.Range("M113").Value = -251.0715
$ws.Range("N113").Value = -7839.5

# Row 136
$ws.Range("H136").Value = 8232.375
$ws.Range("I136").Value = 7829.3335
$ws.Range("J136").Value = 9441.5
$ws.Range("K136").Value = 23488.0005
$ws.Range("L136").Value = 28324.5
$ws.Range("M136").Value = -20938.0005
$ws.Range("N136").Value = -33424.5

$ws = $wb.Worksheets.Item("WVR")

# Row 26
$ws.Range("H26").Value = 2000000
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 2000000
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 2000000
$ws.Range("M26").ClearContents()
$ws.Range("N26").Value = -2000586

# Row 64
$ws.Range("H64").Value = 44749.75
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 44749.75
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 44749.75
$ws.Range("N64").Value = -45245.75

# Row 67
$ws.Range("H67").Value = 44749.75
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 44749.75
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 44749.75
$ws.Range("N67").Value = -46465.75

# Row 96
$ws.Range("H96").Value = 7000
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 7000
$ws.Range("K96").Value = 0
$ws.Range("L96").Value = 7000
$ws.Range("M96").ClearContents()
$ws.Range("N96").Value = -9746

# Row 111
$ws.Range("H111").Value = 0
$ws.Range("I111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("K111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").ClearContents()

# Row 113
$ws.Range("H113").Value = 686.8
$ws.Range("I113").Value = 585.3
$ws.Range("J113").Value = 889.8
$ws.Range("K113").Value = 1755.9
$ws.Range("L113").Value = 2669.4
$ws.Range("M113").Value = 414.1000000000001
$ws.Range("N113").Value = -7009.4
